$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new "Master Content Code" column header in J1 ---
$ws.Range("J1").Value = "Master Content Code"

# Match J1 formatting to the other header cells in row 1 (G1:I1 style)
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New column J values for existing + new data rows ---
$ws.Range("J2").Value = 7

# --- Duplicate row 2 into rows 3 and 4 (same data, plus new Master Content Code) ---
$ws.Range("A2:J2").Copy()
$ws.Range("A3:J3").PasteSpecial(-4104)  # xlPasteAll
$ws.Range("A4:J4").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

# --- Column width for the new column J ---
$ws.Columns("J").ColumnWidth = 21

# --- Sheet view: scroll / selection state ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("C4").Select()
$ws.Rows("4:4").EntireRow.Select()
